$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Insert a new column before column C to hold the Supplier ("Fornecedor") data
$ws.Columns("C:C").Insert()

# Split the store-name merged banner to make room for a CNPJ (store registration) field
$ws.Range("B2:R2").UnMerge()
$ws.Range("B2:O2").Merge()
$ws.Range("P2:S2").Merge()
$ws.Range("P2:S2").HorizontalAlignment = -4108
$ws.Range("P2:S2").VerticalAlignment = -4108
$ws.Range("P2").Value2 = "<CNPJ>"

# New "Fornecedor" (Supplier) header + placeholder column
$ws.Range("C9").Value2 = "Fornecedor"
$ws.Range("C10").Value2 = "<Fornecedor>"

# Refresh the AutoFilter so its range covers the newly inserted column
$ws.AutoFilterMode = $false
$ws.Range("A9:R9").AutoFilter()

# The _FilterDatabase defined name needs to be pointed at the new range as well
$fd = $wb.Names.Item(1)
$fd.RefersTo = "=Plan1!`$A`$9:`$R`$9"

$ws.Range("C11").Select()
